$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "69.727.01"
$ws.Range("E2").Value = "  +0.84%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "2.508.33"
$ws.Range("E3").Value = "  +0.64%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.06%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.51"
$ws.Range("E5").Value = "  +0.26%  "

# --- Row 6 (Solana) ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.64"
$ws.Range("E6").Value = "  +0.74%  "

# --- Row 7 (USDC) ---
$ws.Range("E7").Value = "  +0.00%  "

# --- Row 8 (XRP) ---
$ws.Range("E8").Value = "  +0.24%  "

# --- Row 9 (LidoStakedEther) ---
$ws.Range("D9").Value = "2.506.87"
$ws.Range("E9").Value = "  +0.69%  "

# --- Row 10 (Dogecoin) ---
$ws.Range("E10").Value = "  +4.38%  "

# --- Row 11 (TRON) ---
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.167"
$ws.Range("E11").Value = "  +0.01%  "

# --- Row 12 (Cardano) ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("E12").Value = "  +4.96%  "

# --- Row 13 (Toncoin) ---
$ws.Range("E13").Value = "  +2.66%  "

# --- Row 14 (WrappedliquidstakedEther2.0) ---
$ws.Range("D14").Value = "2.969.40"
$ws.Range("E14").Value = "  +0.82%  "

# --- Row 15 (ShibaInu) ---
$ws.Range("E15").Value = "  +3.22%  "

# --- Row 16 (WrappedBTC) ---
$ws.Range("D16").Value = "69.531.68"
$ws.Range("E16").Value = "  +0.62%  "

# --- Row 17 (Avalanche) ---
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.85"
$ws.Range("E17").Value = "  +1.22%  "

# --- Row 18 (WrappedEther) ---
$ws.Range("D18").Value = "2.502.19"
$ws.Range("E18").Value = "  -0.13%  "

# --- Row 19 (Chainlink) ---
$ws.Range("E19").Value = "  -0.74%  "

# --- Row 20 (Uniswap) ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("E20").Value = "  -2.61%  "

# --- Row 21 (BitcoinCash) ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.36"
$ws.Range("E21").Value = "  +1.14%  "

# --- Row 22 (Polkadot) ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.91"
$ws.Range("E22").Value = "  +0.08%  "

# --- Row 23 (SuiNetwork) ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.95"
$ws.Range("E23").Value = "  +1.10%  "

# --- Row 24 (Dai) ---
$ws.Range("E24").Value = "  +0.04%  "

# --- Row 25 (Litecoin) ---
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.34"
$ws.Range("E25").Value = "  +3.49%  "

# --- Row 26 (NEARProtocol) ---
$ws.Range("E26").Value = "  +0.33%  "

# --- Row 27 (Aptos) ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.82"
$ws.Range("E27").Value = "  -0.22%  "

# --- Row 28 (WrappedeETH) ---
$ws.Range("D28").Value = "2.646.25"
$ws.Range("E28").Value = "  +0.32%  "

# --- Row 29 (Binance-PegBSC-USD) ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"

# --- Row 30 (PEPE) ---
$ws.Range("D30").Value = "0.0₃0893"
$ws.Range("E30").Value = "  +0.51%  "

# --- Row 31 (InternetComputer(DFINITY)) ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.85"
$ws.Range("E31").Value = "  +0.96%  "

# --- Row 32 (Bittensor) ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "457.99"
$ws.Range("E32").Value = "  -1.63%  "

# --- Row 33 (Fetch.AI) ---
$ws.Range("E33").Value = "  -2.52%  "

# --- Row 34 (PancakeSwap) ---
$ws.Range("E34").Value = "  +0.01%  "

# --- Row 35 (FirstDigitalUSD) ---
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.10%  "

# --- Row 36 (Monero) ---
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.76"
$ws.Range("E36").Value = "  +4.09%  "

# --- Row 37 (Kaspa) ---
$ws.Range("E37").Value = "  +1.76%  "

# --- Row 38 (WhiteBITCoin) ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.06"
$ws.Range("E38").Value = "  +0.76%  "

# --- Row 39 (EthereumClassic) ---
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.49"
$ws.Range("E39").Value = "  +1.27%  "

# --- Row 40 (USDe) ---
$ws.Range("E40").Value = "  +0.05%  "

# --- Row 41 (PolygonEcosystemToken) ---
$ws.Range("E41").Value = "  +1.72%  "

# --- Row 42 (RenderToken) ---
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.69"
$ws.Range("E42").Value = "  -0.05%  "

# --- Row 43 (Stacks) ---
$ws.Range("E43").Value = "  +1.68%  "

# --- Row 44 (OKB) ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.06"
$ws.Range("E44").Value = "  +0.02%  "

# --- Row 45 (dogwifhat) ---
$ws.Range("E45").Value = "  -3.12%  "

# --- Row 46 (ImmutableX) ---
$ws.Range("E46").Value = "  -5.72%  "

# --- Row 47 (Aave) ---
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.66"
$ws.Range("E47").Value = "  -0.16%  "

# --- Rows 48 & 49 swap places (ARBITRUM <-> Filecoin), with new volume values ---
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.47"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.520"
$ws.Range("E49").Value = "  -0.34%  "

# --- Row 50 (Cronos) ---
$ws.Range("E50").Value = "  +0.79%  "

# --- Row 51 (Mantle) ---
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.579"
$ws.Range("E51").Value = "  -0.22%  "

